$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D2 ("Research Interests" for Abhishek Patel) still held the generic
# "Interest details..." placeholder - bring it in line with the other rows
# by copying Akshita Kadam's (D3) value *and* formatting (font/alignment)
# down into D2.
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Value = $ws.Range("D3").Value2

# The "Teaching" column (H2:H5) held placeholder text ("Teaching
# details..."/"Teching a"); replace it with the new placeholder text used
# for the message from the dean & principal section.
$ws.Range("H2").Value = "Project details..."
$ws.Range("H3").Value = "Project details..."
$ws.Range("H4").Value = "Project details..."
$ws.Range("H5").Value = "Project details..."

# Leave the selection where the author left it after editing D2:D3.
$ws.Range("D2:D3").Select() | Out-Null
$ws.Range("D3").Activate() | Out-Null
